$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (D in header label) imputed-value swaps on the top block of rows ---
$ws.Cells.Item(2, 5).Value = -7.2
$ws.Cells.Item(3, 5).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(11, 5).Value = -7.9
$ws.Cells.Item(13, 5).Value = ""
$ws.Cells.Item(21, 5).Value = -8.699999999999999
$ws.Cells.Item(25, 5).Value = ""

# --- Remove the "RM 232" row (row 26) and the "SC 92" row (originally row 28,
#     becomes row 27 after the first deletion shifts rows up) ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Fix-ups on the rows that shifted up into their new positions ---
# SC 119 (now row 29): column D (header "C") goes blank
$ws.Cells.Item(29, 4).Value = ""

# SC 232 (now row 33): columns D and E (headers "C" and "D") get real values
$ws.Cells.Item(33, 4).Value = -14.1
$ws.Cells.Item(33, 5).Value = -10.7
